# "add personal info - room info"
# Sheet "빈소2" (room 2) gets its order list replaced:
#   before: 맥주(beer) x3, 치킨(chicken) x1, 대패삼겹살(pork belly) x1  (rows 2-4)
#   after:  치즈김밥(cheese gimbap) x2, 대패삼겹살(pork belly) x1        (rows 2-3)
# Row 1 stays blank (present but empty), row 4 is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("빈소2")

# Wipe out the old order rows (2-4) so stale cells don't linger.
$ws.Range("A2:E4").ClearContents()

# Touch row 1 (no value) so it still materializes as an empty row, matching
# the original layout where data starts on row 2.
$ws.Rows.Item(1).OutlineLevel = 0

# Row 2: 치즈김밥 (cheese gimbap) - 2500 won x 2 = 5000
$ws.Cells.Item(2, 1).Value = "치즈김밥"
$ws.Cells.Item(2, 2).Value = 2500
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = 5000
$ws.Cells.Item(2, 5).Value = "None"

# Row 3: 대패삼겹살 (sliced pork belly) - 2000 won x 1 = 2000
$ws.Cells.Item(3, 1).Value = "대패삼겹살"
$ws.Cells.Item(3, 2).Value = 2000
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 2000
$ws.Cells.Item(3, 5).Value = "None"
